$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "runs", "balls", "fours" for rows 2-4 (Murali Vijay, Chennai Super Kings)
# Leading apostrophe keeps these as text values (matching the original
# numberStoredAsText / t="str" cells) instead of Excel auto-converting
# the numeric-looking strings into real numbers.
$ws.Range("C2").Value = "'21"
$ws.Range("D2").Value = "'21"
$ws.Range("E2").Value = "'3"

$ws.Range("C3").Value = "'1"
$ws.Range("D3").Value = "'7"
$ws.Range("E3").Value = "'0"

$ws.Range("C4").Value = "'10"
$ws.Range("D4").Value = "'15"
$ws.Range("E4").Value = "'1"
